# Generate Report for Archive
# Update localization status for the two files that have moved from
# "Ready for handoff" to "In Translation": 443d0ac9-...md and
# 4718eab0-...md. The third file (a05d68bc-...md) remains "Ready for
# handoff".

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: rows 8 and 9 hold 443d0ac9...md and 4718eab0...md,
# columns B (zh-cn) and C (de-de) carry the per-locale status.
$overview.Range("B8").Value = "In Translation"
$overview.Range("C8").Value = "In Translation"
$overview.Range("B9").Value = "In Translation"
$overview.Range("C9").Value = "In Translation"

# zh-cn / de-de detail sheets: column C is "Status" for rows 8 and 9.
$zhcn.Range("C8").Value = "In Translation"
$zhcn.Range("C9").Value = "In Translation"

$dede.Range("C8").Value = "In Translation"
$dede.Range("C9").Value = "In Translation"
